$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update column F (想去人数) values
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 537
$wsExpo.Range("F5").Value = 271
$wsExpo.Range("F6").Value = 390
$wsExpo.Range("F7").Value = 240
$wsExpo.Range("F8").Value = 2300
$wsExpo.Range("F9").Value = 385
$wsExpo.Range("F10").Value = 5713
$wsExpo.Range("F12").Value = 374

# Sheet "全部类型" (sheet4): update column F (想去人数) values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 537
$wsAll.Range("F6").Value = 271
$wsAll.Range("F7").Value = 390
$wsAll.Range("F8").Value = 240
$wsAll.Range("F11").Value = 2300
$wsAll.Range("F12").Value = 385
$wsAll.Range("F13").Value = 5713
$wsAll.Range("F15").Value = 374
